$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.662.46"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.818.55"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.51"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.566"
$ws.Range("E6").Value = "  +1.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.99"
$ws.Range("E8").Value = "  +8.32%  "

$ws.Range("E9").Value = "  +1.85%  "

$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "2.080.97"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.40"
$ws.Range("E13").Value = "  +3.67%  "

$ws.Range("D14").Value = "1.801.37"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("E15").Value = "  +2.71%  "

$ws.Range("D16").Value = "34.715.54"
$ws.Range("E16").Value = "  +1.26%  "

$ws.Range("E17").Value = "  +3.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.31"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.96"
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.56"
$ws.Range("E21").Value = "  +5.62%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.20"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.64"
$ws.Range("E24").Value = "  +5.71%  "

$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E26").Value = "  +4.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.78"
$ws.Range("E27").Value = "  +2.55%  "

$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.99"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("E31").Value = "  +2.29%  "

$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("E34").Value = "  +2.93%  "

$ws.Range("E35").Value = "  +0.81%  "

$ws.Range("D36").Value = "1.421.57"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.678"
$ws.Range("E37").Value = "  +2.65%  "

$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "86.52"
$ws.Range("E39").Value = "  +5.22%  "

$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("E41").Value = "  +4.39%  "

$ws.Range("E42").Value = "  +4.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.87"
$ws.Range("E44").Value = "  -2.01%  "

$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").Value = "  +2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.13"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("D48").Value = "1.981.14"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.03"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("E51").Value = "  -0.04%  "
